# Term-paper progress: add column K (parity-count helper column) with its
# own conditional formatting, matching the author's next iteration of the
# "курсачтаблица" (term paper table) workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New column K: "IF(J<>0, SUM(B:F), "")" for every data row (2..33)
# ---------------------------------------------------------------------------
$ws.Range("K2").Formula = '=IF(J2<>0, SUM(B2:F2), "")'
$ws.Range("K3:K33").Formula = '=IF(J3<>0, SUM(B3:F3), "")'

# Give column K about the same narrow width as the other helper columns.
$ws.Columns.Item(11).ColumnWidth = 1.6666666666666667

# ---------------------------------------------------------------------------
# 2. Re-prioritise the existing J2:J33 "Equal To" rules so the four new
#    K2:K33 rules can slot in above them (matches the diff's priority/dxfId
#    renumbering: the J rules move from 2/1 down to 6/5).
# ---------------------------------------------------------------------------
$rngJ = $ws.Range("J2:J33")
$jEqual1 = $rngJ.FormatConditions.Item(1)
$jEqual0 = $rngJ.FormatConditions.Item(2)
$jEqual1.Priority = 6
$jEqual0.Priority = 5

# ---------------------------------------------------------------------------
# 3. New K2:K33 "Equal To" conditional formatting, one rule per possible
#    parity-sum value (2,3,4,5), using the standard Excel "Highlight Cells
#    Rules" presets:
#      =2 -> Green Fill with Dark Green Text
#      =3 -> Yellow Fill with Dark Yellow Text
#      =4 -> Light Red Fill with Dark Red Text
#      =5 -> Red Text (no fill)
#    Added in that order, then pinned to priorities 4,3,2,1 so the last
#    rule added (=5) ends up on top, exactly like the committed workbook.
# ---------------------------------------------------------------------------
$rngK = $ws.Range("K2:K33")

$kEqual2 = $rngK.FormatConditions.Add(1, 3, "=2")
$kEqual2.Font.Color = 24832
$kEqual2.Interior.Color = 13561798

$kEqual3 = $rngK.FormatConditions.Add(1, 3, "=3")
$kEqual3.Font.Color = 22428
$kEqual3.Interior.Color = 10284031

$kEqual4 = $rngK.FormatConditions.Add(1, 3, "=4")
$kEqual4.Font.Color = 393372
$kEqual4.Interior.Color = 13551615

$kEqual5 = $rngK.FormatConditions.Add(1, 3, "=5")
$kEqual5.Font.Color = 393372

$kEqual2.Priority = 4
$kEqual3.Priority = 3
$kEqual4.Priority = 2
$kEqual5.Priority = 1

# ---------------------------------------------------------------------------
# 4. Cosmetic: reflect the on-screen selection/scroll state from the author's
#    next editing session (selecting B17:F18 while scrolled down a bit).
# ---------------------------------------------------------------------------
$ws.Range("B17:F18").Select()

Write-Output "done"
